$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: date and invoice number (plain text, not numeric-looking) ---
$ws.Range("F2").Value = "2021-06-06 "
$ws.Range("F3").Value = "INV202106061149"

# --- Bill To details ---
$ws.Range("A9").Value = "Anurag Deo"
$ws.Range("A11").Value = "Bangalore South"

# A13 looks numeric ("990019344") - must stay as text, not be coerced to a number.
# Build it as a formula returning a text string, then flatten to a static value
# via copy / paste-special-values so the style (and hence number format) is left
# completely untouched while the cell keeps a text type.
$ws.Range("A13").Formula = '="990019344"'
$ws.Range("A13").Copy()
$ws.Range("A13").PasteSpecial(-4163)

# --- Salesperson row ---
$ws.Range("B16").Value = "Anurag Deo"

$ws.Range("D16").Formula = '="990019344"'
$ws.Range("D16").Copy()
$ws.Range("D16").PasteSpecial(-4163)

# --- Item row 19 (all numeric-looking values, kept as text) ---
$ws.Range("A19").Formula = '="1"'
$ws.Range("A19").Copy()
$ws.Range("A19").PasteSpecial(-4163)

$ws.Range("B19").Value = "Test Book 9"

$ws.Range("C19").Formula = '="2"'
$ws.Range("C19").Copy()
$ws.Range("C19").PasteSpecial(-4163)

$ws.Range("D19").Formula = '="100"'
$ws.Range("D19").Copy()
$ws.Range("D19").PasteSpecial(-4163)

$ws.Range("E19").Formula = '="0.0"'
$ws.Range("E19").Copy()
$ws.Range("E19").PasteSpecial(-4163)

$ws.Range("F19").Formula = '="200"'
$ws.Range("F19").Copy()
$ws.Range("F19").PasteSpecial(-4163)

# --- Total (already a text cell in the original, keep it text) ---
$ws.Range("F29").Formula = '="200"'
$ws.Range("F29").Copy()
$ws.Range("F29").PasteSpecial(-4163)

$excel.CutCopyMode = 0
